# "Aggiunta elementi mancanti UI" - fill in the missing weekly-report rows
# in the "Formulario Ore" sheet (Dal/Al/Stima/Effettive/Lavoro svolto table).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Week of row 6 ("Requisiti funzionali"): record actual hours worked (D6)
# and extend the work-done note with the next task that was started.
$ws.Range("D6").Value = 16
$ws.Range("E6").Value = "Requisiti funzionali`nInizio sviluppo interfaccia"

# Week of row 7 was missing its actual-hours and work-done entries entirely.
$ws.Range("D7").Value = 16
$ws.Range("E7").Value = "Sviluppo interfaccia`nInizio sviluppo comunicazione Bluetooth"

# Week of row 8 was missing its work-done entry.
$ws.Range("E8").Value = "Sviluppo interfaccia`nGestione permessi dispositivo per Bluetooth"

# Row 4's note gets clarified with the technology that was ultimately chosen.
$ws.Range("E4").Value = "Analisi tecnologie per sviluppo mobile`nSviluppo prototipi con Xamarin e React Native`nScelta finale tecnologia (Xamarin)`n"

# Totals row: sum the "Effettive" (actual hours) column now that it's populated.
$ws.Range("D28").Formula = "=SUM(D3:D27)"

# Leave the selection where the user last worked.
$ws.Range("E8").Select()
